$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. TCDs sheet: trim the shared-formula / lookup block from columns J:N
# ---------------------------------------------------------------------------
$tcds = $wb.Worksheets.Item("TCDs")
$tcds.Activate()

# Row 1 (header "rule-..." formulas) - columns J:N lose their formula but
# keep their explicit style (s="3").
$tcds.Range("J1:N1").ClearContents()

# Rows 2 and 3 (numeric 5 / 7 values) - columns J:N fall back to the column's
# default style, so a full clear (contents + formatting) collapses them back
# to "no override" cells.
$tcds.Range("J2:N3").Clear()

# Row 4 - J:M fall back to the column default style; N4 keeps its own
# style (s="5") so only its contents are cleared.
$tcds.Range("J4:M4").Clear()
$tcds.Range("N4").ClearContents()

# Row 5 keeps its own row style (s="2") across J:N, so only contents clear.
$tcds.Range("J5:N5").ClearContents()

# Row 6 falls back to the column default style again.
$tcds.Range("J6:N6").Clear()

# Restore the frozen-pane scroll position / refresh the active selection.
$win = $excel.ActiveWindow
$pane = $win.Panes.Item(4)
$pane.ScrollRow = 2
$pane.ScrollColumn = 5
$tcds.Range("C16").Select()

# ---------------------------------------------------------------------------
# 2. Hide the supporting lookup sheets - only "TCDs" stays visible.
# ---------------------------------------------------------------------------
$tcdsDict = $wb.Worksheets.Item("TCDs_dict")
$tcdsArch = $wb.Worksheets.Item("TCDs_arch")
$lookups = $wb.Worksheets.Item("Lookups")

# ---------------------------------------------------------------------------
# 3. TCDs_dict - restore its scroll position and refresh its selection.
# ---------------------------------------------------------------------------
$tcdsDict.Activate()
$win = $excel.ActiveWindow
$pane = $win.Panes.Item(4)
$pane.ScrollRow = 27
$pane.ScrollColumn = 5
$tcdsDict.Range("F41:F56").Select()

# ---------------------------------------------------------------------------
# 4. TCDs_arch - restore its scroll position and move the selection.
# ---------------------------------------------------------------------------
$tcdsArch.Activate()
$win = $excel.ActiveWindow
$pane = $win.Panes.Item(4)
$pane.ScrollRow = 22
$pane.ScrollColumn = 7
$tcdsArch.Range("F41:F56").Select()

# ---------------------------------------------------------------------------
# 5. Lookups - move the selection.
# ---------------------------------------------------------------------------
$lookups.Activate()
$lookups.Range("F41:F56").Select()

# Now hide the three helper sheets.
$tcdsDict.Visible = $false
$tcdsArch.Visible = $false
$lookups.Visible = $false

# ---------------------------------------------------------------------------
# 6. Leave "TCDs" as the active / displayed sheet.
# ---------------------------------------------------------------------------
$tcds.Activate()
